$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update SQL query text in shared-string backed cells ---
# Replace the old ".id" join conditions with the new ".study_id" / ".participant_id" ones
# across every cell that contains one of the canned SQL queries.

$replacements = @(
    @('df_participant prt ON std.id = prt."study.id"', 'df_participant prt ON std.study_id = prt."study.study_id"'),
    @('df_diagnoses dgn ON prt.id = dgn."participant.id"', 'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"'),
    @('df_treatments trt ON prt.id = trt."participant.id"', 'df_treatments trt ON prt.participant_id = trt."participant.participant_id"'),
    @('df_treatment_resp trr ON prt.id = trr."participant.id"', 'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"'),
    @('df_survival srv ON prt.id = srv."participant.id"', 'df_survival srv ON prt.participant_id = srv."participant.participant_id"'),
    @('df_reference_files rfs ON std.id = rfs."study.id"', 'df_reference_files rfs ON std.study_id = rfs."study.study_id"')
)

$queryCells = @("B2","C2","B3","B4","B5","B6","B7")

foreach ($addr in $queryCells) {
    $val = $ws.Range($addr).Value2
    if ($val -ne $null) {
        foreach ($pair in $replacements) {
            $val = $val -replace [regex]::Escape($pair[0]), $pair[1]
        }
        $ws.Range($addr).Value2 = $val
    }
}

# --- Sheet view: scroll back to column A instead of B (topLeftCell B6 -> A6) ---
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1

# --- Column C width: drop bestFit, widen to ~68.33203125 ---
$ws.Columns("C").ColumnWidth = 67.5

# --- Row 2 height: 297 -> 336 ---
$ws.Rows("2").RowHeight = 336
